# Apply the tracked-changes edit: move the "Área de Login" slide (sldId 257)
# so it comes after the "Cadastro do fornecedor" slide (sldId 260), i.e. move
# it from position 3 to position 5; then remove the two leftover empty
# placeholder shapes (Título 1 / Subtítulo 2) that remained on the slide
# that is now at position 6 (sldId 261).

$p = $ppt.ActivePresentation

function Get-SlideById($pres, $sldId) {
    for ($i = 1; $i -le $pres.Slides.Count; $i++) {
        $candidate = $pres.Slides.Item($i)
        if ($candidate.SlideID -eq $sldId) {
            return $candidate
        }
    }
    return $null
}

# Move "Área de Login" (sldId 257) to just after "Cadastro do fornecedor"
# (sldId 260), i.e. to slide position 5.
$loginSlide = Get-SlideById $p 257
$loginSlide.MoveTo(5)

# The slide that keeps sldId 261 now ends up at position 6. It was left with
# two empty, content-less placeholder shapes (an unused ctrTitle and
# subTitle) - remove them so the slide has no shapes at all.
$emptySlide = Get-SlideById $p 261
while ($emptySlide.Shapes.Count -gt 0) {
    $emptySlide.Shapes.Item(1).Delete()
}
